$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing REG-TC-02 row (row 5): precondition text now also
# states the page the scenario starts from, and the cell needs to wrap.
$ws.Cells.Item(5, 7).Value = "On registration page and `nemail already exists"
$ws.Cells.Item(5, 7).WrapText = $true

# --- Add the new REG-TC-03 test case as row 6, cloning the formatting
# used by the row above it (row 5) so borders/fonts/alignment all match.
$ws.Range("B5:K5").Copy()
$ws.Range("B6:K6").PasteSpecial(-4122)
$ws.Range("B6:K6").RowHeight = 47.25

$ws.Cells.Item(6, 2).Value  = "REG-03"
$ws.Cells.Item(6, 3).Value  = "REG-TC-03"
$ws.Cells.Item(6, 4).Value  = "Registration Module"
$ws.Cells.Item(6, 5).Value  = "REG-03"
$ws.Cells.Item(6, 6).Value  = "Register with invalid email format"
$ws.Cells.Item(6, 7).Value  = "On registration page"
$ws.Cells.Item(6, 8).Value  = "1. Input invalid email `n2. Input password `n3. Click register"
$ws.Cells.Item(6, 9).Value  = "usermail.com"
$ws.Cells.Item(6, 10).Value = "Error message appears"
$ws.Cells.Item(6, 11).Value = "Medium"

$ws.Hyperlinks.Add($ws.Cells.Item(6, 9), "mailto:user1@gmail.com", "", "", "user1@gmail.com")

$ws.Range("K6").Select()
